$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 667.5
$ws.Range("I19").Value = 696.6667
$ws.Range("J19").Value = 650
$ws.Range("K19").Value = 696.6667
$ws.Range("L19").Value = 650
$ws.Range("M19").Value = -521.6667
$ws.Range("N19").Value = -1000
$ws.Range("H53").Value = 169.3158
$ws.Range("I53").Value = 83.72727
$ws.Range("J53").Value = 287
$ws.Range("K53").Value = 83.72727
$ws.Range("L53").Value = 287
$ws.Range("M53").Value = 553.27273
$ws.Range("N53").Value = -1561
$ws.Range("H62").Value = 3120.625
$ws.Range("I62").Value = 3343.1667
$ws.Range("J62").Value = 2453
$ws.Range("K62").Value = 3343.1667
$ws.Range("L62").Value = 2453
$ws.Range("M62").Value = -2719.1667
$ws.Range("N62").Value = -3701
$ws.Range("H65").Value = 3120.625
$ws.Range("I65").Value = 3343.1667
$ws.Range("J65").Value = 2453
$ws.Range("K65").Value = 16715.8335
$ws.Range("L65").Value = 12265
$ws.Range("M65").Value = -13595.8335
$ws.Range("N65").Value = -18505
$ws.Range("H98").Value = 9000.462
$ws.Range("I98").Value = 5818.727
$ws.Range("J98").Value = 26500
$ws.Range("K98").Value = 5818.727
$ws.Range("L98").Value = 26500
$ws.Range("M98").Value = -4320.727
$ws.Range("N98").Value = -29496
$ws.Range("H100").Value = 3128.2273
$ws.Range("I100").Value = 2860.4167
$ws.Range("J100").Value = 3449.6
$ws.Range("K100").Value = 2860.4167
$ws.Range("L100").Value = 3449.6
$ws.Range("M100").Value = -2319.4167
$ws.Range("N100").Value = -4531.6
$ws.Range("H116").Value = 11767071
$ws.Range("I116").Value = 22224334
$ws.Range("J116").Value = 2648.75
$ws.Range("K116").Value = 22224334
$ws.Range("L116").Value = 2648.75
$ws.Range("M116").Value = -22220892
$ws.Range("N116").Value = -9532.75
$ws.Range("H122").Value = 9000.462
$ws.Range("I122").Value = 5818.727
$ws.Range("J122").Value = 26500
$ws.Range("K122").Value = 17456.181
$ws.Range("L122").Value = 79500
$ws.Range("M122").Value = -15006.181
$ws.Range("N122").Value = -84400
$ws.Range("H138").Value = 1624.5962
$ws.Range("I138").Value = 1232.25
$ws.Range("J138").Value = 2932.4167
$ws.Range("K138").Value = 3696.75
$ws.Range("L138").Value = 8797.250100000001
$ws.Range("M138").Value = 1443.25
$ws.Range("N138").Value = -19077.2501
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17550.5
$ws.Range("I32").Value = 19728.611
$ws.Range("J32").Value = 4481.8335
$ws.Range("K32").Value = 19728.611
$ws.Range("L32").Value = 4481.8335
$ws.Range("M32").Value = -19441.611
$ws.Range("N32").Value = -5055.8335
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2375.875
$ws.Range("I99").Value = 1502.5
$ws.Range("J99").Value = 3249.25
$ws.Range("K99").Value = 1502.5
$ws.Range("L99").Value = 3249.25
$ws.Range("M99").Value = -4.5
$ws.Range("N99").Value = -6245.25
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 109.125
$ws.Range("I7").Value = 81.22221999999999
$ws.Range("J7").Value = 145
$ws.Range("K7").Value = 81.22221999999999
$ws.Range("L7").Value = 145
$ws.Range("M7").Value = 31.77778000000001
$ws.Range("N7").Value = -371
$ws.Range("H122").Value = 3075.2222
$ws.Range("I122").Value = 2468
$ws.Range("J122").Value = 3682.4443
$ws.Range("K122").Value = 7404
$ws.Range("L122").Value = 11047.3329
$ws.Range("M122").Value = -4954
$ws.Range("N122").Value = -15947.3329
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2107.75
$ws.Range("I5").Value = 2107.75
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 6323.25
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -6211.25
$ws.Range("N5").ClearContents()
$ws.Range("H13").Value = 2355.6
$ws.Range("I13").Value = 1360
$ws.Range("J13").Value = 3849
$ws.Range("K13").Value = 4080
$ws.Range("L13").Value = 11547
$ws.Range("M13").Value = -3912
$ws.Range("N13").Value = -11883
$ws.Range("H25").Value = 4442.25
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 4442.25
$ws.Range("K25").Value = 0
$ws.Range("L25").ClearContents()
$ws.Range("M25").Value = 13326.75
$ws.Range("N25").Value = -13664.75
$ws.Range("H30").Value = 4442.25
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 4442.25
$ws.Range("K30").Value = 0
$ws.Range("L30").ClearContents()
$ws.Range("M30").Value = 13326.75
$ws.Range("N30").Value = -13530.75
$ws.Range("H34").Value = 675.65
$ws.Range("I34").Value = 277.69232
$ws.Range("J34").Value = 1414.7142
$ws.Range("K34").Value = 833.07696
$ws.Range("L34").Value = 4244.142599999999
$ws.Range("M34").Value = -749.07696
$ws.Range("N34").Value = -4412.142599999999
$ws.Range("H122").Value = 1166.3334
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 1166.3334
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 10497.0006
$ws.Range("N122").Value = -15397.0006
$ws.Range("H135").Value = 2107.75
$ws.Range("I135").Value = 2107.75
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 18969.75
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -16434.75
$ws.Range("N135").ClearContents()
$ws.Range("H139").Value = 2318.8276
$ws.Range("I139").Value = 1982.2222
$ws.Range("J139").Value = 2869.6365
$ws.Range("K139").Value = 5946.6666
$ws.Range("L139").Value = 8608.9095
$ws.Range("M139").Value = -806.6665999999996
$ws.Range("N139").Value = -18888.9095
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1741
$ws.Range("I107").Value = 2982.2856
$ws.Range("J107").Value = 292.83334
$ws.Range("K107").Value = 2982.2856
$ws.Range("L107").Value = 292.83334
$ws.Range("M107").Value = -1062.2856
$ws.Range("N107").Value = -4132.83334
$ws.Range("H113").Value = 4967.8887
$ws.Range("I113").Value = 3427.75
$ws.Range("J113").Value = 6200
$ws.Range("K113").Value = 3427.75
$ws.Range("L113").Value = 6200
$ws.Range("M113").Value = -1257.75
$ws.Range("N113").Value = -10540
$ws.Range("H136").Value = 43250.2
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 43250.2
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 129750.6
$ws.Range("N136").Value = -134850.6
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2149.8333
$ws.Range("I61").Value = 2426
$ws.Range("J61").Value = 1597.5
$ws.Range("K61").Value = 2426
$ws.Range("L61").Value = 1597.5
$ws.Range("M61").Value = -2224
$ws.Range("N61").Value = -2001.5
$ws.Range("H113").Value = 2149.8333
$ws.Range("I113").Value = 2426
$ws.Range("J113").Value = 1597.5
$ws.Range("K113").Value = 2426
$ws.Range("L113").Value = 1597.5
$ws.Range("M113").Value = -256
$ws.Range("N113").Value = -5937.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 438.52942
$ws.Range("I107").Value = 414.64285
$ws.Range("J107").Value = 550
$ws.Range("K107").Value = 1243.92855
$ws.Range("L107").Value = 1650
$ws.Range("M107").Value = 676.0714499999999
$ws.Range("N107").Value = -5490
$ws.Range("H137").Value = 29761.666
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 29761.666
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 29761.666
$ws.Range("N137").Value = -39961.666
